{"js": "// Remove the empty paragraph, the page-break paragraph, and the\n// \"\u00a9 2020 ... Creative Commons Attribution\" footer paragraph that follow\n// the \"LOB1012: Estat\u00edstica (Requisito fraco)\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"LOB1012: Estat\u00edstica (Requisito fraco)\";\nconst footerText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nconst items = paragraphs.items;\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\n// Only proceed if the footer paragraph we expect to remove (three blocks\n// after the marker) is actually present, so this is a no-op otherwise.\nif (\n  markerIndex !== -1 &&\n  markerIndex + 3 < items.length &&\n  items[markerIndex + 3].text === footerText\n) {\n  // Delete starting from the furthest paragraph so earlier deletions don't\n  // shift the indices of the ones still pending.\n  items[markerIndex + 3].delete();\n  items[markerIndex + 2].delete();\n  items[markerIndex + 1].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the empty paragraph, the page-break paragraph, and the\n# \"\u00a9 2020 ... Creative Commons Attribution\" footer paragraph that follow\n# the \"LOB1012: Estat\u00edstica (Requisito fraco)\" paragraph.\n$d = $word.ActiveDocument\n\n$markerIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*LOB1012*Estat*stica*Requisito fraco*\") {\n        $markerIndex = $i\n        break\n    }\n}\n\n$footerIndex = $markerIndex + 3\n$footerPresent = $false\nif ($markerIndex -ne -1 -and $footerIndex -le $d.Paragraphs.Count) {\n    $footerPresent = $d.Paragraphs.Item($footerIndex).Range.Text -like \"*2020*Contact*luizeleno*usp.br*Creative Commons Attribution*\"\n}\n\nif ($markerIndex -ne -1 -and $footerPresent) {\n    # Delete the three paragraphs right after the marker, starting from the\n    # furthest one so earlier deletions don't shift later indices.\n    for ($offset = 3; $offset -ge 1; $offset--) {\n        $idx = $markerIndex + $offset\n        if ($idx -le $d.Paragraphs.Count) {\n            $d.Paragraphs.Item($idx).Range.Delete()\n        }\n    }\n}\n"}
